$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 11 (pushes existing rows 11+ down by 3,
# carrying row heights/styles along - mirrors native Excel "Insert Rows")
$ws.Rows("11:13").Insert()

# Row 11: new rule r7.0 (Crew Ratings intro)
$ws.Range("A11").Value = 'r7.0'
$ws.Range("B11").Value = '<Bold>r7.0 Crew Ratings</Bold> 
<LineBreak/><LineBreak/>
Each crew member is given a numerical rating of his skill at his position. The rating range from 1 (poorest) to 10 (best). 
<LineBreak/><LineBreak/>
The higher the crew member''s rating, the more successful he will be trying to spot enemy units, hit enemy units with his weapons, repair malfunctioning guns, etc. 
<LineBreak/><LineBreak/>The crew ratings are recorded on the After Action Report (AAR) <InlineUIContainer><Button Content=''r2.4'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>.
<LineBreak/><LineBreak/>
<InlineUIContainer><Button Content=''r7.1'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer> Ratings for New Men<LineBreak/>
<InlineUIContainer><Button Content=''r7.2'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer> Rating Improvements<LineBreak/>'
$ws.Rows(11).RowHeight = 142.65

# Row 12: new rule r7.1 (Ratings for New Men) - label first, content follows below
$ws.Range("A12").Value = 'r7.1'

# Row 13: new rule r7.2 (Ratings Improvement) - label first, content follows below
$ws.Range("A13").Value = 'r7.2'

$ws.Range("B12").Value = '<Bold>r7.1 Ratings for New Men</Bold> 
<LineBreak/><LineBreak/>
Determine the rating for a new crew member by rolling 1D/2 rounded up. For example, a roll of five results in a rating of 3.  
<LineBreak/><LineBreak/>When rolling for more than one new crewman, assign the rating after all have been rolled for. 
<LineBreak/><LineBreak/>
For example, if rolling for an entire crew, roll five ratings and then assign to men as you wish. It is suggested that the higher ratings be allocated in this order: Commander, Gunner, Driver, Loader, and Assistant Driver.'
$ws.Rows(12).RowHeight = 99.85

$ws.Range("B13").Value = '<Bold>r7.2 Ratings Improvement</Bold> 
<LineBreak/><LineBreak/>
At the end of the day of action, check for each surviving crew member to see if his rating improves. Roll 1D separately for each. If the number rolled is higher than the current rating, the rating is increased by one. 
<LineBreak/><LineBreak/>
For example, your gunner has a rating of 5 and you roll a 9. The gunner''s rating increases to 6.
<LineBreak/><LineBreak/>
Crew ratings cannot be greater than 10 and never decrease. Ratings cannot improve if you avoid action for the day per <InlineUIContainer><Button Content=''r20.1'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>.
<LineBreak/><LineBreak/>
Ratings may also improve during refitting periods per <InlineUIContainer><Button Content=''r27.1'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>.'
$ws.Rows(13).RowHeight = 156.9

# Refresh the persisted sort-state metadata to the shifted range (A65:B240 -> A68:B243)
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A68:A243"))
$sortObj.SetRange($ws.Range("A68:B243"))
$sortObj.Apply()
